# Update automatico via Actualizar 02-15-2021 12-08-48
# Shifts the "last updated" timestamps stored in column D down one
# generation and stamps the most recent batch of rows with the newest
# refresh time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newest = 44242.50600596993
$middle = 44242.48478574074
$oldest = 44242.46356857639

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $middle
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldest
}
